# Update the CondFile entry for the 3-back (3) condition: the referenced
# file was renamed from "threeback(3)_2.xlsx" to "threeback(3)_v2.xlsx".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Conditions/Version_2/threeback(3)_v2.xlsx"
